# ICTPMG613 E-Pert Chart: rename existing sheet to "Chart" and add two new
# sheets, "k-Milestones" and "k-Dependences", each holding a small table.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original (only) sheet to "Chart" -----------------------
$chart = $wb.Worksheets.Item(1)
$chart.Name = "Chart"

# --- 2. Add "k-Milestones" right after "Chart" -----------------------------
$milestones = $wb.Worksheets.Add($null, $chart)
$milestones.Name = "k-Milestones"

# --- 3. Add "k-Dependences" right after "k-Milestones" ---------------------
$dependences = $wb.Worksheets.Add($null, $milestones)
$dependences.Name = "k-Dependences"

# ============================================================================
# k-Milestones sheet
# ============================================================================

$milestoneNames = @(
    "Project Charters Approved",
    "Kickoff Meetings Conducted",
    "Project Plans Approved",
    "Cloud Provider Contract Signed",
    "Cloud Infrastructure Configured",
    "Databases Successfully Migrated",
    "Core Applications Deployed to Cloud",
    "Devices Distributed and Configured",
    "Training Completed",
    "First Sprint Completed & Reviewed",
    "Website Deployed to Production",
    "Formal Project Closure"
)

$dueDates = @(45874, 45877, 45877, 45881, 45889, 45889, 45897, 45901, 45917, 45898, 45901, 45929)

# Header row
$milestones.Range("A1").Value = "No."
$milestones.Range("B1").Value = "Milestone"
$milestones.Range("C1").Value = "Due Date"

$hdr = $milestones.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1

# Data rows (2 .. 13)
$milestones.Range("A2").Value = 1
$milestones.Range("A3").Formula = "=A2+1"
$milestones.Range("A4:A13").Formula = "=A3+1"

for ($i = 0; $i -lt $milestoneNames.Length; $i++) {
    $r = $i + 2
    $milestones.Range("B" + $r).Value = $milestoneNames[$i]
    $milestones.Range("C" + $r).Value = $dueDates[$i]
}

$body = $milestones.Range("A2:C13")
$body.Borders.LineStyle = 1
$body.VerticalAlignment = -4160

$dates = $milestones.Range("C2:C13")
$dates.NumberFormat = "ddd\ dd/mmm/yy"

$milestones.Columns(2).ColumnWidth = 29.5
$milestones.Columns(3).ColumnWidth = 13.61

$milestones.Range("A2:A13").Select()

# ============================================================================
# k-Dependences sheet
# ============================================================================

$dependenceNames = @(
    "Gain Formal Acceptance of Project Charter",
    "Gain Formal Acceptance of Scope Statement",
    "Gain Formal Acceptance of Project Plans",
    "Formal budget obtaining",
    "Select Cloud Provider and Sign Contract",
    "Migrate Databases to a Cloud platform",
    "Develop and Deploy Core Applications to Cloud platform",
    "Perform System-Wide Testing",
    "Select work devices Provider and Sign Contract",
    "Configure Work Devices",
    "Gain Formal Acceptance of website functionalities",
    "Conduct Final Testing: IT infrastructure and website, and remote connection"
)

$dependences.Range("A1").Value = "No."
$dependences.Range("B1").Value = "Depedency"

$hdr2 = $dependences.Range("A1:B1")
$hdr2.Font.Bold = $true
$hdr2.Borders.LineStyle = 1

$dependences.Range("A2").Value = 1
$dependences.Range("A3").Formula = "=A2+1"
$dependences.Range("A4:A13").Formula = "=A3+1"

for ($i = 0; $i -lt $dependenceNames.Length; $i++) {
    $r = $i + 2
    $dependences.Range("B" + $r).Value = $dependenceNames[$i]
}

$body2 = $dependences.Range("A2:B13")
$body2.Borders.LineStyle = 1
$body2.VerticalAlignment = -4160

$dependences.Columns(1).ColumnWidth = 3.05
$dependences.Columns(2).ColumnWidth = 62.39

$dependences.Range("A3").Select()

# ============================================================================
# Restore Chart's original view (remove top-left/selection quirks introduced
# by activating the other sheets) and keep k-Dependences as the active tab.
# ============================================================================

$chart.Range("A60").Select()
$chart.Range("G4").Select()

$dependences.Activate()
$dependences.Range("A3").Select()
